# Apply the "ADC calculations x1000 instead of x1000000" change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new row at row 16. This pushes the existing rows 16-25
#    down to 17-26 and Excel automatically re-points every formula
#    that referenced those rows (B20*1000 -> B21*1000, B22/B19 ->
#    B23/B20, etc.) which is exactly the row shift seen in the diff.
# ------------------------------------------------------------------
$ws.Rows("16:16").Insert()

# ------------------------------------------------------------------
# 2) Populate the brand new row 16 ("div/mA" / "Div/mA"): values in
#    row 15 (Div/A) divided by 1000.
# ------------------------------------------------------------------
$ws.Range("A16").Value = "div/mA"
$ws.Range("B16").Formula = "=B15/1000"
$ws.Range("C16").Formula = "=C15/1000"
$ws.Range("D16").Formula = "=D15/1000"
$ws.Range("E16").Value = "Div/mA"
$ws.Range("F16").Formula = "=F15/1000"
$ws.Range("G16").Formula = "=G15/1000"
$ws.Range("H16").Formula = "=H15/1000"
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""

# ------------------------------------------------------------------
# 3) Row 21 (previously row 20, the raw ADC "Count" values row) gets
#    new values: C/D/H change (1200 -> 4095 style full-scale values).
# ------------------------------------------------------------------
$ws.Range("B21").Value = 60
$ws.Range("C21").Value = 3650
$ws.Range("D21").Value = 4095
$ws.Range("F21").Value = 60
$ws.Range("G21").Value = 3388
$ws.Range("H21").Value = 4095

# ------------------------------------------------------------------
# 4) Row 22 ("Count * 1000") now also has F/G/H populated (mirroring
#    B/C/D) and no longer has a stray E22 value.
# ------------------------------------------------------------------
$ws.Range("E22").Value = ""
$ws.Range("F22").Formula = "=F21*1000"
$ws.Range("G22").Formula = "=G21*1000"
$ws.Range("H22").Formula = "=H21*1000"

# Rows 23 (Count * 1000000) and 24 (mA) formulas already point at the
# correct (shifted) source rows automatically after the row insert
# above, so nothing else needs to change there.

# ------------------------------------------------------------------
# 5) Row 26 (previously a blank spacer row) becomes a raw copy of the
#    "Count * 1000" row (row 22) values, feeding the overflow checks.
# ------------------------------------------------------------------
$ws.Range("B26").Formula = "=B22"
$ws.Range("C26").Formula = "=C22"
$ws.Range("D26").Formula = "=D22"
$ws.Range("E26").Value = ""
$ws.Range("F26").Formula = "=F22"
$ws.Range("G26").Formula = "=G22"
$ws.Range("H26").Formula = "=H22"
$ws.Range("I26").Value = ""
$ws.Range("J26").Value = ""

# ------------------------------------------------------------------
# 6) New rows 27-29: overflow checks for int16_t / int32_t / int64_t.
# ------------------------------------------------------------------
$ws.Range("A27").Value = "int16_t"
$ws.Range("B27").Formula = '=IF(B26>(2^15),"OVF","OK")'
$ws.Range("C27").Formula = '=IF(C26>(2^15),"OVF","OK")'
$ws.Range("D27").Formula = '=IF(D26>(2^15),"OVF","OK")'
$ws.Range("E27").Formula = "=2^15"
$ws.Range("F27").Formula = '=IF(F26>(2^15),"OVF","OK")'
$ws.Range("G27").Formula = '=IF(G26>(2^15),"OVF","OK")'
$ws.Range("H27").Formula = '=IF(H26>(2^15),"OVF","OK")'

$ws.Range("A28").Value = "int32_t"
$ws.Range("B28").Formula = '=IF(B26>(2^31),"OVF","OK")'
$ws.Range("C28").Formula = '=IF(C26>(2^31),"OVF","OK")'
$ws.Range("D28").Formula = '=IF(D26>(2^31),"OVF","OK")'
$ws.Range("E28").Formula = "=2^31"
$ws.Range("F28").Formula = '=IF(F26>(2^31),"OVF","OK")'
$ws.Range("G28").Formula = '=IF(G26>(2^31),"OVF","OK")'
$ws.Range("H28").Formula = '=IF(H26>(2^31),"OVF","OK")'

$ws.Range("A29").Value = "int64_t"
$ws.Range("B29").Formula = '=IF(B26>(2^63),"OVF","OK")'
$ws.Range("C29").Formula = '=IF(C26>(2^63),"OVF","OK")'
$ws.Range("D29").Formula = '=IF(D26>(2^63),"OVF","OK")'
$ws.Range("E29").Formula = "=2^63"
$ws.Range("F29").Formula = '=IF(F26>(2^63),"OVF","OK")'
$ws.Range("G29").Formula = '=IF(G26>(2^63),"OVF","OK")'
$ws.Range("H29").Formula = '=IF(H26>(2^63),"OVF","OK")'

# ------------------------------------------------------------------
# 7) Sheet-level cosmetics that accompanied the edit in the diff.
# ------------------------------------------------------------------
$ws.Range("C16").Select()
$ws.Columns("A:A").ColumnWidth = 14.88

Write-Host "Done."
